$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.375.20'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.625.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.485'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.850.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.623.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.387.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.64'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0518'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.98%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.36%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.159.60'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.81%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.804'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.64%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.783'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.763.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +9.52%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.410'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('E50').Style = 'Normal'

